$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.807494584125266
$ws.Range("C2").Value = 0.7062961560329768
$ws.Range("D2").Value = 0.07600777982561624
$ws.Range("E2").Value = 0.01739878643989634
$ws.Range("G2").Value = 0.002637015469981312
$ws.Range("L2").Value = 0.3498838686687975
$ws.Range("N2").Value = 3.80890238265917
$ws.Range("B3").Value = 5.562066824044734
$ws.Range("C3").Value = 0.6448770912427335
$ws.Range("D3").Value = 0.06929318894320602
$ws.Range("E3").Value = 0.01702221701724671
$ws.Range("G3").Value = 0.0026469214176816
$ws.Range("L3").Value = 0.3386790311696757
$ws.Range("N3").Value = 3.723955635705124
$ws.Range("B4").Value = 5.415710663929815
$ws.Range("C4").Value = 0.6076830267203377
$ws.Range("D4").Value = 0.06522548548886675
$ws.Range("E4").Value = 0.01678808838690937
$ws.Range("G4").Value = 0.002653304275436961
$ws.Range("L4").Value = 0.3320503415256155
$ws.Range("N4").Value = 3.672218710653596
$ws.Range("B5").Value = 5.357142615910789
$ws.Range("C5").Value = 0.5926520124864396
$ws.Range("D5").Value = 0.06358125906319856
$ws.Range("E5").Value = 0.0166919173938096
$ws.Range("G5").Value = 0.002655981281615113
$ws.Range("L5").Value = 0.3294114350602797
$ws.Range("N5").Value = 3.65123677361214
$ws.Range("B6").Value = 5.347481758104891
$ws.Range("C6").Value = 0.5901636034387252
$ws.Range("D6").Value = 0.06330903246255559
$ws.Range("E6").Value = 0.01667590143106334
$ws.Range("G6").Value = 0.002656430393612655
$ws.Range("L6").Value = 0.328976988467943
$ws.Range("N6").Value = 3.647758725239783
$ws.Range("B7").Value = 5.414916471800609
$ws.Range("C7").Value = 0.6074798091291882
$ws.Range("D7").Value = 0.06520325723214171
$ws.Range("E7").Value = 0.01678679451340948
$ws.Range("G7").Value = 0.002653340070574985
$ws.Range("L7").Value = 0.3320145009580671
$ws.Range("N7").Value = 3.671935336756832
$ws.Range("B8").Value = 5.721960157234776
$ws.Range("C8").Value = 0.6850088497817524
$ws.Range("D8").Value = 0.07368087774500509
$ws.Range("E8").Value = 0.01726953226483374
$ws.Range("G8").Value = 0.00264036887654831
$ws.Range("L8").Value = 0.3459678222941562
$ws.Range("N8").Value = 3.779522664776721
$ws.Range("B9").Value = 6.359327038464244
$ws.Range("C9").Value = 0.8413524699259938
$ws.Range("D9").Value = 0.09076427000651677
$ws.Range("E9").Value = 0.01819442570346919
$ws.Range("G9").Value = 0.002617300719547358
$ws.Range("L9").Value = 0.3753619521917102
$ws.Range("N9").Value = 3.994055139814066
$ws.Range("B10").Value = 6.850334258304656
$ws.Range("C10").Value = 0.9591435908566268
$ws.Range("D10").Value = 0.1036270563958368
$ws.Range("E10").Value = 0.01886265003681142
$ws.Range("G10").Value = 0.002601773013958246
$ws.Range("L10").Value = 0.3982544958580974
$ws.Range("N10").Value = 4.154155626409192
$ws.Range("B11").Value = 7.078923933356236
$ws.Range("C11").Value = 1.013432339139001
$ws.Range("D11").Value = 0.1095535314171201
$ws.Range("E11").Value = 0.01916463100381272
$ws.Range("G11").Value = 0.002595012442081855
$ws.Range("L11").Value = 0.4089636625303399
$ws.Range("N11").Value = 4.227599514158783
$ws.Range("B12").Value = 7.16625941638506
$ws.Range("C12").Value = 1.034096895793539
$ws.Range("D12").Value = 0.1118091182229932
$ws.Range("E12").Value = 0.01927873228002142
$ws.Range("G12").Value = 0.002592495570919111
$ws.Range("L12").Value = 0.4130624651136827
$ws.Range("N12").Value = 4.255504596899073
$ws.Range("B13").Value = 7.147415400655404
$ws.Range("C13").Value = 1.029641588725383
$ws.Range("D13").Value = 0.1113228230123298
$ws.Range("E13").Value = 0.01925416918629708
$ws.Range("G13").Value = 0.00259303570809154
$ws.Range("L13").Value = 0.4121777645559064
$ws.Range("N13").Value = 4.249490495416467
$ws.Range("B14").Value = 7.086093446298264
$ws.Range("C14").Value = 1.015130258106751
$ws.Range("D14").Value = 0.1097388689994858
$ws.Range("E14").Value = 0.01917402303860971
$ws.Range("G14").Value = 0.002594804513667308
$ws.Range("L14").Value = 0.4092999962518036
$ws.Range("N14").Value = 4.229893378088377
$ws.Range("B15").Value = 7.048633372842232
$ws.Range("C15").Value = 1.006255687015198
$ws.Range("D15").Value = 0.1087701475276361
$ws.Range("E15").Value = 0.01912489939663953
$ws.Range("G15").Value = 0.002595893574391717
$ws.Range("L15").Value = 0.407542972435607
$ws.Range("N15").Value = 4.217901909768557
$ws.Range("B16").Value = 6.835502430445558
$ws.Range("C16").Value = 0.9556103214083009
$ws.Range("D16").Value = 0.1032413062556543
$ws.Range("E16").Value = 0.01884287763355008
$ws.Range("G16").Value = 0.002602220900081219
$ws.Range("L16").Value = 0.397560653874848
$ws.Range("N16").Value = 4.149368686449151
$ws.Range("B17").Value = 6.706107784943697
$ws.Range("C17").Value = 0.9247252514278443
$ws.Range("D17").Value = 0.0998691740299904
$ws.Range("E17").Value = 0.01866937929537027
$ws.Range("G17").Value = 0.002606179868431199
$ws.Range("L17").Value = 0.3915131561434748
$ws.Range("N17").Value = 4.107486589398945
$ws.Range("B18").Value = 6.632174152206858
$ws.Range("C18").Value = 0.9070270109802436
$ws.Range("D18").Value = 0.0979366500476857
$ws.Range("E18").Value = 0.01856939746730024
$ws.Range("G18").Value = 0.002608485509978054
$ws.Range("L18").Value = 0.3880625212545254
$ws.Range("N18").Value = 4.083454526232003
$ws.Range("B19").Value = 6.607225128779305
$ws.Range("C19").Value = 0.9010458761534323
$ws.Range("D19").Value = 0.09728352262142437
$ws.Range("E19").Value = 0.01853551162907774
$ws.Range("G19").Value = 0.002609271074863502
$ws.Range("L19").Value = 0.3868989269950305
$ws.Range("N19").Value = 4.075327384154662
$ws.Range("B20").Value = 6.719831120237131
$ws.Range("C20").Value = 0.9280061482218684
$ws.Range("D20").Value = 0.1002274109525274
$ws.Range("E20").Value = 0.01868786792209054
$ws.Range("G20").Value = 0.002605755477666675
$ws.Range("L20").Value = 0.3921540458008508
$ws.Range("N20").Value = 4.111939030446223
$ws.Range("B21").Value = 7.104084020508481
$ws.Range("C21").Value = 1.019389652672601
$ws.Range("D21").Value = 0.1102038018897957
$ws.Range("E21").Value = 0.0191975704710412
$ws.Range("G21").Value = 0.002594283802867501
$ws.Range("L21").Value = 0.4101440784133956
$ws.Range("N21").Value = 4.235646946941074
$ws.Range("B22").Value = 7.359734842183343
$ws.Range("C22").Value = 1.079737848658908
$ws.Range("D22").Value = 0.1167904306106493
$ws.Range("E22").Value = 0.01952923650006433
$ws.Range("G22").Value = 0.002587038091210852
$ws.Range("L22").Value = 0.4221554989326819
$ws.Range("N22").Value = 4.317044606242803
$ws.Range("B23").Value = 7.222868300862842
$ws.Range("C23").Value = 1.047470041606346
$ws.Range("D23").Value = 0.11326875202775
$ws.Range("E23").Value = 0.01935234135543062
$ws.Range("G23").Value = 0.002590882356724222
$ws.Range("L23").Value = 0.4157212009307614
$ws.Range("N23").Value = 4.273549284373246
$ws.Range("B24").Value = 6.71362537830862
$ws.Range("C24").Value = 0.9265226745243922
$ws.Range("D24").Value = 0.1000654328958177
$ws.Range("E24").Value = 0.01867950994953471
$ws.Range("G24").Value = 0.002605947252452842
$ws.Range("L24").Value = 0.3918642182050718
$ws.Range("N24").Value = 4.109925937370747
$ws.Range("B25").Value = 6.183001964177492
$ws.Range("C25").Value = 0.7985656693661554
$ws.Range("D25").Value = 0.08609039293601484
$ws.Range("E25").Value = 0.01794634616326274
$ws.Range("G25").Value = 0.002623290125233198
$ws.Range("L25").Value = 0.3671866527062377
$ws.Range("N25").Value = 3.935605542588775
